$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style/border/number-format) of the last existing data
# row (385) down across the new rows (386-464) before writing values, so the
# new rows inherit the same cell styles (date format in column A, etc.)
$srcRow = $ws.Range("A385:D385")
$dstRange = $ws.Range("A386:D464")
$srcRow.Copy()
$dstRange.PasteSpecial(-4122)

# New daily data rows (date serial, nuovi pos., somma mobile 7gg.,
# somma mobile 7gg. per 100mila abitanti) updated through 2021-12-08.
$newData = @(
    @(386, 44460, 2, 4, 163.1986944104447),
    @(387, 44461, 0, 4, 163.1986944104447),
    @(388, 44462, 0, 4, 163.1986944104447),
    @(389, 44463, 1, 5, 203.9983680130559),
    @(390, 44464, 0, 5, 203.9983680130559),
    @(391, 44465, 0, 3, 122.3990208078335),
    @(392, 44466, 1, 4, 163.1986944104447),
    @(393, 44467, 0, 2, 81.59934720522236),
    @(394, 44468, 0, 2, 81.59934720522236),
    @(395, 44469, 1, 3, 122.3990208078335),
    @(396, 44470, 0, 2, 81.59934720522236),
    @(397, 44471, 0, 2, 81.59934720522236),
    @(398, 44472, 0, 2, 81.59934720522236),
    @(399, 44473, 0, 1, 40.79967360261118),
    @(400, 44474, 0, 1, 40.79967360261118),
    @(401, 44475, 0, 1, 40.79967360261118),
    @(402, 44476, 0, 0, 0.0),
    @(403, 44477, 0, 0, 0.0),
    @(404, 44478, 0, 0, 0.0),
    @(405, 44479, 0, 0, 0.0),
    @(406, 44480, 0, 0, 0.0),
    @(407, 44481, 0, 0, 0.0),
    @(408, 44482, 0, 0, 0.0),
    @(409, 44483, 0, 0, 0.0),
    @(410, 44484, 0, 0, 0.0),
    @(411, 44485, 0, 0, 0.0),
    @(412, 44486, 0, 0, 0.0),
    @(413, 44487, 0, 0, 0.0),
    @(414, 44488, 0, 0, 0.0),
    @(415, 44489, 0, 0, 0.0),
    @(416, 44490, 0, 0, 0.0),
    @(417, 44491, 0, 0, 0.0),
    @(418, 44492, 0, 0, 0.0),
    @(419, 44493, 0, 0, 0.0),
    @(420, 44494, 0, 0, 0.0),
    @(421, 44495, 0, 0, 0.0),
    @(422, 44496, 0, 0, 0.0),
    @(423, 44497, 1, 1, 40.79967360261118),
    @(424, 44498, 0, 1, 40.79967360261118),
    @(425, 44499, 0, 1, 40.79967360261118),
    @(426, 44500, 0, 1, 40.79967360261118),
    @(427, 44501, 2, 3, 122.3990208078335),
    @(428, 44502, 0, 3, 122.3990208078335),
    @(429, 44503, 0, 3, 122.3990208078335),
    @(430, 44504, 0, 2, 81.59934720522236),
    @(431, 44505, 0, 2, 81.59934720522236),
    @(432, 44506, 0, 2, 81.59934720522236),
    @(433, 44507, 1, 3, 122.3990208078335),
    @(434, 44508, 0, 1, 40.79967360261118),
    @(435, 44509, 0, 1, 40.79967360261118),
    @(436, 44510, 0, 1, 40.79967360261118),
    @(437, 44511, 3, 4, 163.1986944104447),
    @(438, 44512, 0, 4, 163.1986944104447),
    @(439, 44513, 0, 4, 163.1986944104447),
    @(440, 44514, 0, 3, 122.3990208078335),
    @(441, 44515, 0, 3, 122.3990208078335),
    @(442, 44516, 1, 4, 163.1986944104447),
    @(443, 44517, 0, 4, 163.1986944104447),
    @(444, 44518, 0, 1, 40.79967360261118),
    @(445, 44519, 0, 1, 40.79967360261118),
    @(446, 44520, 0, 1, 40.79967360261118),
    @(447, 44521, 1, 2, 81.59934720522236),
    @(448, 44522, 0, 2, 81.59934720522236),
    @(449, 44523, 1, 2, 81.59934720522236),
    @(450, 44524, 3, 5, 203.9983680130559),
    @(451, 44525, 1, 6, 244.798041615667),
    @(452, 44526, 0, 6, 244.798041615667),
    @(453, 44527, 3, 9, 367.1970624235006),
    @(454, 44528, 4, 12, 489.5960832313341),
    @(455, 44529, 1, 13, 530.3957568339454),
    @(456, 44530, 0, 12, 489.5960832313341),
    @(457, 44531, 0, 9, 367.1970624235006),
    @(458, 44532, 2, 10, 407.9967360261118),
    @(459, 44533, 3, 13, 530.3957568339454),
    @(460, 44534, 1, 11, 448.796409628723),
    @(461, 44535, 0, 7, 285.5977152182783),
    @(462, 44536, 2, 8, 326.3973888208894),
    @(463, 44537, 1, 9, 367.1970624235006),
    @(464, 44538, 1, 10, 407.9967360261118)
)

foreach ($row in $newData) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

Write-Host "Added rows 386-464"
